$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value2 = "Datos actualizados a 15 de Abril de 2020 a las 15:22"

$rows = 213
$cols = 8
$data = New-Object 'object[,]' $rows,$cols

$data[0,0] = "Estados Unidos"
$data[0,1] = 614246
$data[0,2] = 360
$data[0,3] = 38820
$data[0,4] = 549362
$data[0,5] = 13473
$data[0,6] = 17
$data[0,7] = 26064
$data[1,0] = "España"
$data[1,1] = 177633
$data[1,2] = 3573
$data[1,3] = 70853
$data[1,4] = 88201
$data[1,5] = 7371
$data[1,6] = 324
$data[1,7] = 18579
$data[2,0] = "Italia"
$data[2,1] = 162488
$data[2,2] = 0
$data[2,3] = 37130
$data[2,4] = 104291
$data[2,5] = 3186
$data[2,6] = 0
$data[2,7] = 21067
$data[3,0] = "Francia"
$data[3,1] = 143303
$data[3,2] = 0
$data[3,3] = 28805
$data[3,4] = 98769
$data[3,5] = 6730
$data[3,6] = 0
$data[3,7] = 15729
$data[4,0] = "Alemania"
$data[4,1] = 132210
$data[4,2] = 0
$data[4,3] = 72600
$data[4,4] = 56115
$data[4,5] = 4288
$data[4,6] = 0
$data[4,7] = 3495
$data[5,0] = "Reino Unido"
$data[5,1] = 93873
$data[5,2] = 0
$data[5,3] = 0
$data[5,4] = 81422
$data[5,5] = 1559
$data[5,6] = 0
$data[5,7] = 12107
$data[6,0] = "China"
$data[6,1] = 82295
$data[6,2] = 46
$data[6,3] = 77816
$data[6,4] = 1137
$data[6,5] = 113
$data[6,6] = 1
$data[6,7] = 3342
$data[7,0] = "Iran"
$data[7,1] = 76389
$data[7,2] = 1512
$data[7,3] = 49933
$data[7,4] = 21679
$data[7,5] = 3643
$data[7,6] = 94
$data[7,7] = 4777
$data[8,0] = "Turquia"
$data[8,1] = 65111
$data[8,2] = 0
$data[8,3] = 4799
$data[8,4] = 58909
$data[8,5] = 1809
$data[8,6] = 0
$data[8,7] = 1403
$data[9,0] = "Belgica"
$data[9,1] = 33573
$data[9,2] = 2454
$data[9,3] = 7107
$data[9,4] = 22026
$data[9,5] = 1204
$data[9,6] = 283
$data[9,7] = 4440
$data[10,0] = "Paises Bajos"
$data[10,1] = 28153
$data[10,2] = 734
$data[10,3] = 250
$data[10,4] = 24769
$data[10,5] = 1358
$data[10,6] = 189
$data[10,7] = 3134
$data[11,0] = "Canada"
$data[11,1] = 27063
$data[11,2] = 0
$data[11,3] = 8235
$data[11,4] = 17925
$data[11,5] = 557
$data[11,6] = 0
$data[11,7] = 903
$data[12,0] = "Suiza"
$data[12,1] = 26336
$data[12,2] = 400
$data[12,3] = 14700
$data[12,4] = 10415
$data[12,5] = 386
$data[12,6] = 47
$data[12,7] = 1221
$data[13,0] = "Brasil"
$data[13,1] = 25758
$data[13,2] = 496
$data[13,3] = 14026
$data[13,4] = 10175
$data[13,5] = 296
$data[13,6] = 25
$data[13,7] = 1557
$data[14,0] = "Rusia"
$data[14,1] = 24490
$data[14,2] = 3388
$data[14,3] = 1986
$data[14,4] = 22306
$data[14,5] = 8
$data[14,6] = 28
$data[14,7] = 198
$data[15,0] = "Portugal"
$data[15,1] = 18091
$data[15,2] = 643
$data[15,3] = 383
$data[15,4] = 17109
$data[15,5] = 208
$data[15,6] = 32
$data[15,7] = 599
$data[16,0] = "Austria"
$data[16,1] = 14321
$data[16,2] = 95
$data[16,3] = 8098
$data[16,4] = 5830
$data[16,5] = 232
$data[16,6] = 9
$data[16,7] = 393
$data[17,0] = "Israel"
$data[17,1] = 12200
$data[17,2] = 154
$data[17,3] = 2309
$data[17,4] = 9765
$data[17,5] = 176
$data[17,6] = 3
$data[17,7] = 126
$data[18,0] = "Suecia"
$data[18,1] = 11927
$data[18,2] = 482
$data[18,3] = 381
$data[18,4] = 10343
$data[18,5] = 954
$data[18,6] = 170
$data[18,7] = 1203
$data[19,0] = "India"
$data[19,1] = 11555
$data[19,2] = 68
$data[19,3] = 1362
$data[19,4] = 9797
$data[19,5] = 0
$data[19,6] = 3
$data[19,7] = 396
$data[20,0] = "Irlanda"
$data[20,1] = 11479
$data[20,2] = 0
$data[20,3] = 77
$data[20,4] = 10996
$data[20,5] = 194
$data[20,6] = 0
$data[20,7] = 406
$data[21,0] = "Corea del Sur"
$data[21,1] = 10591
$data[21,2] = 27
$data[21,3] = 7616
$data[21,4] = 2750
$data[21,5] = 55
$data[21,6] = 3
$data[21,7] = 225
$data[22,0] = "Peru"
$data[22,1] = 10303
$data[22,2] = 0
$data[22,3] = 2869
$data[22,4] = 7204
$data[22,5] = 132
$data[22,6] = 0
$data[22,7] = 230
$data[23,0] = "Japon"
$data[23,1] = 8100
$data[23,2] = 215
$data[23,3] = 853
$data[23,4] = 7101
$data[23,5] = 152
$data[23,6] = 0
$data[23,7] = 146
$data[24,0] = "Chile"
$data[24,1] = 7917
$data[24,2] = 0
$data[24,3] = 2646
$data[24,4] = 5179
$data[24,5] = 387
$data[24,6] = 0
$data[24,7] = 92
$data[25,0] = "Ecuador"
$data[25,1] = 7603
$data[25,2] = 0
$data[25,3] = 696
$data[25,4] = 6538
$data[25,5] = 129
$data[25,6] = 0
$data[25,7] = 369
$data[26,0] = "Polonia"
$data[26,1] = 7408
$data[26,2] = 206
$data[26,3] = 668
$data[26,4] = 6472
$data[26,5] = 160
$data[26,6] = 5
$data[26,7] = 268
$data[27,0] = "Rumania"
$data[27,1] = 7216
$data[27,2] = 337
$data[27,3] = 1217
$data[27,4] = 5637
$data[27,5] = 245
$data[27,6] = 11
$data[27,7] = 362
$data[28,0] = "Noruega"
$data[28,1] = 6686
$data[28,2] = 63
$data[28,3] = 32
$data[28,4] = 6512
$data[28,5] = 59
$data[28,6] = 3
$data[28,7] = 142
$data[29,0] = "Dinamarca"
$data[29,1] = 6681
$data[29,2] = 170
$data[29,3] = 2515
$data[29,4] = 3857
$data[29,5] = 100
$data[29,6] = 10
$data[29,7] = 309
$data[30,0] = "Australia"
$data[30,1] = 6447
$data[30,2] = 47
$data[30,3] = 3686
$data[30,4] = 2698
$data[30,5] = 76
$data[30,6] = 2
$data[30,7] = 63
$data[31,0] = "Chequia"
$data[31,1] = 6151
$data[31,2] = 40
$data[31,3] = 676
$data[31,4] = 5312
$data[31,5] = 84
$data[31,6] = 2
$data[31,7] = 163
$data[32,0] = "Pakistan"
$data[32,1] = 5988
$data[32,2] = 151
$data[32,3] = 1446
$data[32,4] = 4435
$data[32,5] = 46
$data[32,6] = 11
$data[32,7] = 107
$data[33,0] = "Arabia Saudita"
$data[33,1] = 5862
$data[33,2] = 493
$data[33,3] = 931
$data[33,4] = 4852
$data[33,5] = 59
$data[33,6] = 6
$data[33,7] = 79
$data[34,0] = "Filipinas"
$data[34,1] = 5453
$data[34,2] = 230
$data[34,3] = 353
$data[34,4] = 4751
$data[34,5] = 1
$data[34,6] = 14
$data[34,7] = 349
$data[35,0] = "Mexico"
$data[35,1] = 5399
$data[35,2] = 385
$data[35,3] = 2125
$data[35,4] = 2868
$data[35,5] = 207
$data[35,6] = 74
$data[35,7] = 406
$data[36,0] = "Indonesia"
$data[36,1] = 5136
$data[36,2] = 297
$data[36,3] = 446
$data[36,4] = 4221
$data[36,5] = 0
$data[36,6] = 10
$data[36,7] = 469
$data[37,0] = "Malasia"
$data[37,1] = 5072
$data[37,2] = 85
$data[37,3] = 2647
$data[37,4] = 2342
$data[37,5] = 56
$data[37,6] = 1
$data[37,7] = 83
$data[38,0] = "Emiratos Arabes Unidos"
$data[38,1] = 4933
$data[38,2] = 0
$data[38,3] = 933
$data[38,4] = 3972
$data[38,5] = 1
$data[38,6] = 0
$data[38,7] = 28
$data[39,0] = "Serbia"
$data[39,1] = 4873
$data[39,2] = 408
$data[39,3] = 400
$data[39,4] = 4374
$data[39,5] = 131
$data[39,6] = 5
$data[39,7] = 99
$data[40,0] = "Ucrania"
$data[40,1] = 3764
$data[40,2] = 392
$data[40,3] = 143
$data[40,4] = 3513
$data[40,5] = 45
$data[40,6] = 10
$data[40,7] = 108
$data[41,0] = "Bielorrusia"
$data[41,1] = 3728
$data[41,2] = 447
$data[41,3] = 203
$data[41,4] = 3489
$data[41,5] = 68
$data[41,6] = 3
$data[41,7] = 36
$data[42,0] = "Catar"
$data[42,1] = 3711
$data[42,2] = 283
$data[42,3] = 406
$data[42,4] = 3298
$data[42,5] = 37
$data[42,6] = 0
$data[42,7] = 7
$data[43,0] = "Panama"
$data[43,1] = 3574
$data[43,2] = 0
$data[43,3] = 72
$data[43,4] = 3407
$data[43,5] = 106
$data[43,6] = 0
$data[43,7] = 95
$data[44,0] = "Luxemburgo"
$data[44,1] = 3307
$data[44,2] = 0
$data[44,3] = 500
$data[44,4] = 2740
$data[44,5] = 30
$data[44,6] = 0
$data[44,7] = 67
$data[45,0] = "Republica Dominicana"
$data[45,1] = 3286
$data[45,2] = 0
$data[45,3] = 162
$data[45,4] = 2941
$data[45,5] = 143
$data[45,6] = 0
$data[45,7] = 183
$data[46,0] = "Singapur"
$data[46,1] = 3252
$data[46,2] = 0
$data[46,3] = 611
$data[46,4] = 2631
$data[46,5] = 29
$data[46,6] = 0
$data[46,7] = 10
$data[47,0] = "Finlandia"
$data[47,1] = 3237
$data[47,2] = 76
$data[47,3] = 300
$data[47,4] = 2873
$data[47,5] = 75
$data[47,6] = 0
$data[47,7] = 64
$data[48,0] = "Colombia"
$data[48,1] = 2979
$data[48,2] = 0
$data[48,3] = 354
$data[48,4] = 2498
$data[48,5] = 106
$data[48,6] = 0
$data[48,7] = 127
$data[49,0] = "Tailandia"
$data[49,1] = 2643
$data[49,2] = 30
$data[49,3] = 1497
$data[49,4] = 1103
$data[49,5] = 61
$data[49,6] = 2
$data[49,7] = 43
$data[50,0] = "Argentina"
$data[50,1] = 2443
$data[50,2] = 166
$data[50,3] = 559
$data[50,4] = 1776
$data[50,5] = 83
$data[50,6] = 6
$data[50,7] = 108
$data[51,0] = "Sudafrica"
$data[51,1] = 2415
$data[51,2] = 0
$data[51,3] = 410
$data[51,4] = 1978
$data[51,5] = 7
$data[51,6] = 0
$data[51,7] = 27
$data[52,0] = "Egipto"
$data[52,1] = 2350
$data[52,2] = 0
$data[52,3] = 589
$data[52,4] = 1583
$data[52,5] = 0
$data[52,6] = 0
$data[52,7] = 178
$data[53,0] = "Grecia"
$data[53,1] = 2170
$data[53,2] = 0
$data[53,3] = 269
$data[53,4] = 1800
$data[53,5] = 76
$data[53,6] = 0
$data[53,7] = 101
$data[54,0] = "Argelia"
$data[54,1] = 2070
$data[54,2] = 0
$data[54,3] = 691
$data[54,4] = 1053
$data[54,5] = 60
$data[54,6] = 0
$data[54,7] = 326
$data[55,0] = "Marruecos"
$data[55,1] = 1988
$data[55,2] = 100
$data[55,3] = 218
$data[55,4] = 1643
$data[55,5] = 1
$data[55,6] = 1
$data[55,7] = 127
$data[56,0] = "Moldavia"
$data[56,1] = 1934
$data[56,2] = 0
$data[56,3] = 171
$data[56,4] = 1720
$data[56,5] = 80
$data[56,6] = 2
$data[56,7] = 43
$data[57,0] = "Croacia"
$data[57,1] = 1741
$data[57,2] = 37
$data[57,3] = 473
$data[57,4] = 1234
$data[57,5] = 31
$data[57,6] = 3
$data[57,7] = 34
$data[58,0] = "Islandia"
$data[58,1] = 1720
$data[58,2] = 0
$data[58,3] = 989
$data[58,4] = 723
$data[58,5] = 8
$data[58,6] = 0
$data[58,7] = 8
$data[59,0] = "Barein"
$data[59,1] = 1671
$data[59,2] = 143
$data[59,3] = 663
$data[59,4] = 1001
$data[59,5] = 3
$data[59,6] = 0
$data[59,7] = 7
$data[60,0] = "Hungria"
$data[60,1] = 1579
$data[60,2] = 67
$data[60,3] = 192
$data[60,4] = 1253
$data[60,5] = 58
$data[60,6] = 12
$data[60,7] = 134
$data[61,0] = "Kuwait"
$data[61,1] = 1405
$data[61,2] = 50
$data[61,3] = 206
$data[61,4] = 1196
$data[61,5] = 31
$data[61,6] = 0
$data[61,7] = 3
$data[62,0] = "Estonia"
$data[62,1] = 1400
$data[62,2] = 27
$data[62,3] = 117
$data[62,4] = 1248
$data[62,5] = 10
$data[62,6] = 4
$data[62,7] = 35
$data[63,0] = "Irak"
$data[63,1] = 1400
$data[63,2] = 0
$data[63,3] = 766
$data[63,4] = 556
$data[63,5] = 0
$data[63,6] = 0
$data[63,7] = 78
$data[64,0] = "Nueva Zelanda"
$data[64,1] = 1386
$data[64,2] = 20
$data[64,3] = 728
$data[64,4] = 649
$data[64,5] = 3
$data[64,6] = 0
$data[64,7] = 9
$data[65,0] = "Kazajistan"
$data[65,1] = 1290
$data[65,2] = 58
$data[65,3] = 240
$data[65,4] = 1034
$data[65,5] = 20
$data[65,6] = 2
$data[65,7] = 16
$data[66,0] = "Uzbekistan"
$data[66,1] = 1275
$data[66,2] = 110
$data[66,3] = 107
$data[66,4] = 1164
$data[66,5] = 8
$data[66,6] = 0
$data[66,7] = 4
$data[67,0] = "Eslovenia"
$data[67,1] = 1248
$data[67,2] = 28
$data[67,3] = 165
$data[67,4] = 1022
$data[67,5] = 34
$data[67,6] = 5
$data[67,7] = 61
$data[68,0] = "Banglades"
$data[68,1] = 1231
$data[68,2] = 219
$data[68,3] = 49
$data[68,4] = 1132
$data[68,5] = 1
$data[68,6] = 4
$data[68,7] = 50
$data[69,0] = "Azerbaiyan"
$data[69,1] = 1197
$data[69,2] = 0
$data[69,3] = 351
$data[69,4] = 833
$data[69,5] = 25
$data[69,6] = 0
$data[69,7] = 13
$data[70,0] = "Armenia"
$data[70,1] = 1111
$data[70,2] = 44
$data[70,3] = 297
$data[70,4] = 797
$data[70,5] = 30
$data[70,6] = 1
$data[70,7] = 17
$data[71,0] = "Bosnia y Herzegovina"
$data[71,1] = 1110
$data[71,2] = 27
$data[71,3] = 253
$data[71,4] = 816
$data[71,5] = 4
$data[71,6] = 1
$data[71,7] = 41
$data[72,0] = "Lituania"
$data[72,1] = 1091
$data[72,2] = 21
$data[72,3] = 138
$data[72,4] = 924
$data[72,5] = 14
$data[72,6] = 0
$data[72,7] = 29
$data[73,0] = "Hong Kong"
$data[73,1] = 1017
$data[73,2] = 4
$data[73,3] = 459
$data[73,4] = 554
$data[73,5] = 10
$data[73,6] = 0
$data[73,7] = 4
$data[74,0] = "Republica de Macedonia"
$data[74,1] = 974
$data[74,2] = 66
$data[74,3] = 98
$data[74,4] = 831
$data[74,5] = 15
$data[74,6] = 1
$data[74,7] = 45
$data[75,0] = "Oman"
$data[75,1] = 910
$data[75,2] = 97
$data[75,3] = 131
$data[75,4] = 775
$data[75,5] = 3
$data[75,6] = 0
$data[75,7] = 4
$data[76,0] = "Eslovaquia"
$data[76,1] = 863
$data[76,2] = 28
$data[76,3] = 151
$data[76,4] = 706
$data[76,5] = 5
$data[76,6] = 4
$data[76,7] = 6
$data[77,0] = "Camerun"
$data[77,1] = 848
$data[77,2] = 0
$data[77,3] = 165
$data[77,4] = 666
$data[77,5] = 0
$data[77,6] = 3
$data[77,7] = 17
$data[78,0] = "Afganistan"
$data[78,1] = 784
$data[78,2] = 70
$data[78,3] = 43
$data[78,4] = 716
$data[78,5] = 0
$data[78,6] = 2
$data[78,7] = 25
$data[79,0] = "Cuba"
$data[79,1] = 766
$data[79,2] = 0
$data[79,3] = 132
$data[79,4] = 613
$data[79,5] = 9
$data[79,6] = 0
$data[79,7] = 21
$data[80,0] = "Tunez"
$data[80,1] = 747
$data[80,2] = 0
$data[80,3] = 43
$data[80,4] = 670
$data[80,5] = 89
$data[80,6] = 0
$data[80,7] = 34
$data[81,0] = "Bulgaria"
$data[81,1] = 735
$data[81,2] = 22
$data[81,3] = 105
$data[81,4] = 594
$data[81,5] = 29
$data[81,6] = 1
$data[81,7] = 36
$data[82,0] = "Crucero"
$data[82,1] = 712
$data[82,2] = 0
$data[82,3] = 639
$data[82,4] = 61
$data[82,5] = 7
$data[82,6] = 0
$data[82,7] = 12
$data[83,0] = "Republica de Chipre"
$data[83,1] = 695
$data[83,2] = 0
$data[83,3] = 65
$data[83,4] = 618
$data[83,5] = 8
$data[83,6] = 0
$data[83,7] = 12
$data[84,0] = "Letonia"
$data[84,1] = 666
$data[84,2] = 9
$data[84,3] = 44
$data[84,4] = 617
$data[84,5] = 3
$data[84,6] = 0
$data[84,7] = 5
$data[85,0] = "Principado de Andorra"
$data[85,1] = 659
$data[85,2] = 0
$data[85,3] = 128
$data[85,4] = 500
$data[85,5] = 17
$data[85,6] = 0
$data[85,7] = 31
$data[86,0] = "Libano"
$data[86,1] = 658
$data[86,2] = 17
$data[86,3] = 81
$data[86,4] = 556
$data[86,5] = 33
$data[86,6] = 0
$data[86,7] = 21
$data[87,0] = "Costa de Marfil"
$data[87,1] = 638
$data[87,2] = 0
$data[87,3] = 114
$data[87,4] = 518
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 6
$data[88,0] = "Ghana"
$data[88,1] = 636
$data[88,2] = 0
$data[88,3] = 17
$data[88,4] = 611
$data[88,5] = 2
$data[88,6] = 0
$data[88,7] = 8
$data[89,0] = "Costa Rica"
$data[89,1] = 618
$data[89,2] = 0
$data[89,3] = 66
$data[89,4] = 549
$data[89,5] = 14
$data[89,6] = 0
$data[89,7] = 3
$data[90,0] = "Niger"
$data[90,1] = 570
$data[90,2] = 0
$data[90,3] = 90
$data[90,4] = 466
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 14
$data[91,0] = "Burkina Faso"
$data[91,1] = 528
$data[91,2] = 0
$data[91,3] = 177
$data[91,4] = 321
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 30
$data[92,0] = "Albania"
$data[92,1] = 494
$data[92,2] = 19
$data[92,3] = 251
$data[92,4] = 218
$data[92,5] = 5
$data[92,6] = 1
$data[92,7] = 25
$data[93,0] = "Uruguay"
$data[93,1] = 492
$data[93,2] = 0
$data[93,3] = 260
$data[93,4] = 224
$data[93,5] = 14
$data[93,6] = 0
$data[93,7] = 8
$data[94,0] = "Kirguistan"
$data[94,1] = 449
$data[94,2] = 19
$data[94,3] = 78
$data[94,4] = 366
$data[94,5] = 5
$data[94,6] = 0
$data[94,7] = 5
$data[95,0] = "Honduras"
$data[95,1] = 419
$data[95,2] = 12
$data[95,3] = 9
$data[95,4] = 379
$data[95,5] = 10
$data[95,6] = 5
$data[95,7] = 31
$data[96,0] = "Malta"
$data[96,1] = 399
$data[96,2] = 6
$data[96,3] = 44
$data[96,4] = 352
$data[96,5] = 4
$data[96,6] = 0
$data[96,7] = 3
$data[97,0] = "Bolivia"
$data[97,1] = 397
$data[97,2] = 43
$data[97,3] = 7
$data[97,4] = 362
$data[97,5] = 3
$data[97,6] = 0
$data[97,7] = 28
$data[98,0] = "Jordania"
$data[98,1] = 397
$data[98,2] = 0
$data[98,3] = 235
$data[98,4] = 155
$data[98,5] = 5
$data[98,6] = 0
$data[98,7] = 7
$data[99,0] = "Taiwan"
$data[99,1] = 395
$data[99,2] = 2
$data[99,3] = 137
$data[99,4] = 252
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 6
$data[100,0] = "Reunion"
$data[100,1] = 391
$data[100,2] = 0
$data[100,3] = 40
$data[100,4] = 351
$data[100,5] = 3
$data[100,6] = 0
$data[100,7] = 0
$data[101,0] = "Nigeria"
$data[101,1] = 373
$data[101,2] = 0
$data[101,3] = 99
$data[101,4] = 263
$data[101,5] = 2
$data[101,6] = 0
$data[101,7] = 11
$data[102,0] = "San Marino"
$data[102,1] = 372
$data[102,2] = 0
$data[102,3] = 53
$data[102,4] = 283
$data[102,5] = 15
$data[102,6] = 0
$data[102,7] = 36
$data[103,0] = "Guinea"
$data[103,1] = 363
$data[103,2] = 0
$data[103,3] = 31
$data[103,4] = 332
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 0
$data[104,0] = "Republica de Yibuti"
$data[104,1] = 363
$data[104,2] = 0
$data[104,3] = 53
$data[104,4] = 308
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 2
$data[105,0] = "Mauricio"
$data[105,1] = 324
$data[105,2] = 0
$data[105,3] = 51
$data[105,4] = 264
$data[105,5] = 3
$data[105,6] = 0
$data[105,7] = 9
$data[106,0] = "Senegal"
$data[106,1] = 314
$data[106,2] = 15
$data[106,3] = 190
$data[106,4] = 122
$data[106,5] = 1
$data[106,6] = 0
$data[106,7] = 2
$data[107,0] = "Estado de Palestina"
$data[107,1] = 308
$data[107,2] = 0
$data[107,3] = 62
$data[107,4] = 244
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 2
$data[108,0] = "Georgia"
$data[108,1] = 306
$data[108,2] = 6
$data[108,3] = 69
$data[108,4] = 234
$data[108,5] = 6
$data[108,6] = 0
$data[108,7] = 3
$data[109,0] = "Montenegro"
$data[109,1] = 288
$data[109,2] = 5
$data[109,3] = 46
$data[109,4] = 238
$data[109,5] = 7
$data[109,6] = 0
$data[109,7] = 4
$data[110,0] = "Vietnam"
$data[110,1] = 267
$data[110,2] = 1
$data[110,3] = 171
$data[110,4] = 96
$data[110,5] = 8
$data[110,6] = 0
$data[110,7] = 0
$data[111,0] = "Isla de Man"
$data[111,1] = 256
$data[111,2] = 2
$data[111,3] = 151
$data[111,4] = 101
$data[111,5] = 13
$data[111,6] = 2
$data[111,7] = 4
$data[112,0] = "Consejo Danes para los Refugiados"
$data[112,1] = 241
$data[112,2] = 0
$data[112,3] = 20
$data[112,4] = 201
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 20
$data[113,0] = "Sri Lanka"
$data[113,1] = 235
$data[113,2] = 2
$data[113,3] = 63
$data[113,4] = 165
$data[113,5] = 1
$data[113,6] = 0
$data[113,7] = 7
$data[114,0] = "Kenia"
$data[114,1] = 225
$data[114,2] = 9
$data[114,3] = 53
$data[114,4] = 162
$data[114,5] = 2
$data[114,6] = 1
$data[114,7] = 10
$data[115,0] = "Mayotte"
$data[115,1] = 217
$data[115,2] = 0
$data[115,3] = 69
$data[115,4] = 145
$data[115,5] = 3
$data[115,6] = 0
$data[115,7] = 3
$data[116,0] = "Venezuela"
$data[116,1] = 197
$data[116,2] = 8
$data[116,3] = 111
$data[116,4] = 77
$data[116,5] = 6
$data[116,6] = 0
$data[116,7] = 9
$data[117,0] = "Islas Feroe"
$data[117,1] = 184
$data[117,2] = 0
$data[117,3] = 166
$data[117,4] = 18
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 0
$data[118,0] = "Guatemala"
$data[118,1] = 180
$data[118,2] = 13
$data[118,3] = 19
$data[118,4] = 156
$data[118,5] = 3
$data[118,6] = 0
$data[118,7] = 5
$data[119,0] = "Paraguay"
$data[119,1] = 161
$data[119,2] = 2
$data[119,3] = 23
$data[119,4] = 130
$data[119,5] = 1
$data[119,6] = 1
$data[119,7] = 8
$data[120,0] = "El Salvador"
$data[120,1] = 159
$data[120,2] = 10
$data[120,3] = 30
$data[120,4] = 123
$data[120,5] = 2
$data[120,6] = 0
$data[120,7] = 6
$data[121,0] = "Martinica"
$data[121,1] = 158
$data[121,2] = 1
$data[121,3] = 73
$data[121,4] = 77
$data[121,5] = 17
$data[121,6] = 2
$data[121,7] = 8
$data[122,0] = "Mali"
$data[122,1] = 148
$data[122,2] = 4
$data[122,3] = 34
$data[122,4] = 101
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 13
$data[123,0] = "Guadalupe"
$data[123,1] = 145
$data[123,2] = 0
$data[123,3] = 67
$data[123,4] = 70
$data[123,5] = 13
$data[123,6] = 0
$data[123,7] = 8
$data[124,0] = "Brunei"
$data[124,1] = 136
$data[124,2] = 0
$data[124,3] = 108
$data[124,4] = 27
$data[124,5] = 2
$data[124,6] = 0
$data[124,7] = 1
$data[125,0] = "Ruanda"
$data[125,1] = 134
$data[125,2] = 0
$data[125,3] = 49
$data[125,4] = 85
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = "Gibraltar"
$data[126,1] = 129
$data[126,2] = 0
$data[126,3] = 104
$data[126,4] = 25
$data[126,5] = 1
$data[126,6] = 0
$data[126,7] = 0
$data[127,0] = "Camboya"
$data[127,1] = 122
$data[127,2] = 0
$data[127,3] = 96
$data[127,4] = 26
$data[127,5] = 1
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = "Trinidad yTobago"
$data[128,1] = 113
$data[128,2] = 0
$data[128,3] = 17
$data[128,4] = 88
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 8
$data[129,0] = "Madagascar"
$data[129,1] = 108
$data[129,2] = 0
$data[129,3] = 23
$data[129,4] = 85
$data[129,5] = 1
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = "Jamaica"
$data[130,1] = 105
$data[130,2] = 0
$data[130,3] = 21
$data[130,4] = 79
$data[130,5] = 0
$data[130,6] = 1
$data[130,7] = 5
$data[131,0] = "Monaco"
$data[131,1] = 93
$data[131,2] = 0
$data[131,3] = 6
$data[131,4] = 86
$data[131,5] = 5
$data[131,6] = 0
$data[131,7] = 1
$data[132,0] = "Aruba"
$data[132,1] = 92
$data[132,2] = 0
$data[132,3] = 32
$data[132,4] = 60
$data[132,5] = 1
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = "Guayana Francesa"
$data[133,1] = 86
$data[133,2] = 0
$data[133,3] = 51
$data[133,4] = 35
$data[133,5] = 1
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = "Etiopia"
$data[134,1] = 85
$data[134,2] = 3
$data[134,3] = 15
$data[134,4] = 67
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 3
$data[135,0] = "Togo"
$data[135,1] = 81
$data[135,2] = 4
$data[135,3] = 35
$data[135,4] = 43
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 3
$data[136,0] = "Gabon"
$data[136,1] = 80
$data[136,2] = 23
$data[136,3] = 4
$data[136,4] = 75
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 1
$data[137,0] = "Liechtenstein"
$data[137,1] = 79
$data[137,2] = 0
$data[137,3] = 55
$data[137,4] = 23
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 1
$data[138,0] = "Birmania"
$data[138,1] = 74
$data[138,2] = 11
$data[138,3] = 2
$data[138,4] = 68
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 4
$data[139,0] = "Congo"
$data[139,1] = 74
$data[139,2] = 0
$data[139,3] = 10
$data[139,4] = 59
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 5
$data[140,0] = "Barbados"
$data[140,1] = 73
$data[140,2] = 0
$data[140,3] = 15
$data[140,4] = 53
$data[140,5] = 4
$data[140,6] = 0
$data[140,7] = 5
$data[141,0] = "Somalia"
$data[141,1] = 60
$data[141,2] = 0
$data[141,3] = 2
$data[141,4] = 56
$data[141,5] = 2
$data[141,6] = 0
$data[141,7] = 2
$data[142,0] = "Tanzania"
$data[142,1] = 59
$data[142,2] = 6
$data[142,3] = 7
$data[142,4] = 49
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 3
$data[143,0] = "Liberia"
$data[143,1] = 59
$data[143,2] = 0
$data[143,3] = 4
$data[143,4] = 49
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 6
$data[144,0] = "Bermudas"
$data[144,1] = 57
$data[144,2] = 0
$data[144,3] = 30
$data[144,4] = 22
$data[144,5] = 3
$data[144,6] = 0
$data[144,7] = 5
$data[145,0] = "Polinesia Francesa"
$data[145,1] = 55
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 55
$data[145,5] = 1
$data[145,6] = 0
$data[145,7] = 0
$data[146,0] = "Uganda"
$data[146,1] = 55
$data[146,2] = 0
$data[146,3] = 12
$data[146,4] = 43
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = "Islas Caimanes"
$data[147,1] = 54
$data[147,2] = 0
$data[147,3] = 6
$data[147,4] = 47
$data[147,5] = 3
$data[147,6] = 0
$data[147,7] = 1
$data[148,0] = "San Martin (Parte Holandesa)"
$data[148,1] = 52
$data[148,2] = 0
$data[148,3] = 5
$data[148,4] = 38
$data[148,5] = 2
$data[148,6] = 0
$data[148,7] = 9
$data[149,0] = "Bahamas"
$data[149,1] = 49
$data[149,2] = 0
$data[149,3] = 6
$data[149,4] = 35
$data[149,5] = 1
$data[149,6] = 0
$data[149,7] = 8
$data[150,0] = "Zambia"
$data[150,1] = 48
$data[150,2] = 3
$data[150,3] = 30
$data[150,4] = 16
$data[150,5] = 1
$data[150,6] = 0
$data[150,7] = 2
$data[151,0] = "Guyana"
$data[151,1] = 47
$data[151,2] = 0
$data[151,3] = 8
$data[151,4] = 33
$data[151,5] = 5
$data[151,6] = 0
$data[151,7] = 6
$data[152,0] = "Cabo Verde"
$data[152,1] = 45
$data[152,2] = 34
$data[152,3] = 1
$data[152,4] = 43
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 1
$data[153,0] = "Macao"
$data[153,1] = 45
$data[153,2] = 0
$data[153,3] = 10
$data[153,4] = 35
$data[153,5] = 1
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = "Guinea-Bisau"
$data[154,1] = 43
$data[154,2] = 0
$data[154,3] = 0
$data[154,4] = 43
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = "Guinea Ecuatorial"
$data[155,1] = 41
$data[155,2] = 0
$data[155,3] = 4
$data[155,4] = 37
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = "Haiti"
$data[156,1] = 40
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 37
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 3
$data[157,0] = "Puerto Rico"
$data[157,1] = 39
$data[157,2] = 0
$data[157,3] = 1
$data[157,4] = 36
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 2
$data[158,0] = "Eritrea"
$data[158,1] = 35
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 35
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = "Libia"
$data[159,1] = 35
$data[159,2] = 0
$data[159,3] = 9
$data[159,4] = 25
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 1
$data[160,0] = "Benin"
$data[160,1] = 35
$data[160,2] = 0
$data[160,3] = 18
$data[160,4] = 16
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 1
$data[161,0] = "Guam"
$data[161,1] = 32
$data[161,2] = 0
$data[161,3] = 0
$data[161,4] = 31
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 1
$data[162,0] = "Sudan"
$data[162,1] = 32
$data[162,2] = 0
$data[162,3] = 4
$data[162,4] = 23
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 5
$data[163,0] = "San Martin (Parte Francesa)"
$data[163,1] = 32
$data[163,2] = 0
$data[163,3] = 11
$data[163,4] = 19
$data[163,5] = 5
$data[163,6] = 0
$data[163,7] = 2
$data[164,0] = "Mongolia"
$data[164,1] = 30
$data[164,2] = 0
$data[164,3] = 5
$data[164,4] = 25
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = "Siria"
$data[165,1] = 29
$data[165,2] = 0
$data[165,3] = 5
$data[165,4] = 22
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 2
$data[166,0] = "Mozambique"
$data[166,1] = 28
$data[166,2] = 0
$data[166,3] = 2
$data[166,4] = 26
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = "Republica del Chad"
$data[167,1] = 23
$data[167,2] = 0
$data[167,3] = 2
$data[167,4] = 21
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = "Antigua y Barbuda"
$data[168,1] = 23
$data[168,2] = 0
$data[168,3] = 3
$data[168,4] = 18
$data[168,5] = 1
$data[168,6] = 0
$data[168,7] = 2
$data[169,0] = "Maldivas"
$data[169,1] = 21
$data[169,2] = 1
$data[169,3] = 16
$data[169,4] = 5
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = "Laos"
$data[170,1] = 19
$data[170,2] = 0
$data[170,3] = 1
$data[170,4] = 18
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = "Angola"
$data[171,1] = 19
$data[171,2] = 0
$data[171,3] = 5
$data[171,4] = 12
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 2
$data[172,0] = "Nueva Caledonia"
$data[172,1] = 18
$data[172,2] = 0
$data[172,3] = 1
$data[172,4] = 17
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = "Belice"
$data[173,1] = 18
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 16
$data[173,5] = 1
$data[173,6] = 0
$data[173,7] = 2
$data[174,0] = "Zimbabue"
$data[174,1] = 18
$data[174,2] = 1
$data[174,3] = 1
$data[174,4] = 14
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 3
$data[175,0] = "Islas Virgenes de los Estados Unidos"
$data[175,1] = 17
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 17
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = "Fiyi"
$data[176,1] = 16
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 16
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = "Nepal"
$data[177,1] = 16
$data[177,2] = 0
$data[177,3] = 1
$data[177,4] = 15
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = "Malaui"
$data[178,1] = 16
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 14
$data[178,5] = 1
$data[178,6] = 0
$data[178,7] = 2
$data[179,0] = "Namibia"
$data[179,1] = 16
$data[179,2] = 0
$data[179,3] = 3
$data[179,4] = 13
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = "Dominica"
$data[180,1] = 16
$data[180,2] = 0
$data[180,3] = 8
$data[180,4] = 8
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = "Suazilandia"
$data[181,1] = 15
$data[181,2] = 0
$data[181,3] = 8
$data[181,4] = 7
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = "Santa Lucia"
$data[182,1] = 15
$data[182,2] = 0
$data[182,3] = 11
$data[182,4] = 4
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = "San Cristobal y Nieves"
$data[183,1] = 14
$data[183,2] = 0
$data[183,3] = 0
$data[183,4] = 14
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = "Granada"
$data[184,1] = 14
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 14
$data[184,5] = 2
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = "Curazao"
$data[185,1] = 14
$data[185,2] = 0
$data[185,3] = 10
$data[185,4] = 3
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 1
$data[186,0] = "Sierra Leona"
$data[186,1] = 13
$data[186,2] = 2
$data[186,3] = 0
$data[186,4] = 13
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = "Botsuana"
$data[187,1] = 13
$data[187,2] = 0
$data[187,3] = 0
$data[187,4] = 12
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 1
$data[188,0] = "San Vicente y las Granadinas"
$data[188,1] = 12
$data[188,2] = 0
$data[188,3] = 1
$data[188,4] = 11
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 0
$data[189,0] = "Seychelles"
$data[189,1] = 11
$data[189,2] = 0
$data[189,3] = 0
$data[189,4] = 11
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0
$data[190,0] = "Montserrat"
$data[190,1] = 11
$data[190,2] = 0
$data[190,3] = 1
$data[190,4] = 10
$data[190,5] = 1
$data[190,6] = 0
$data[190,7] = 0
$data[191,0] = "Islas Malvinas"
$data[191,1] = 11
$data[191,2] = 0
$data[191,3] = 1
$data[191,4] = 10
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = "Republica de Africa Central"
$data[192,1] = 11
$data[192,2] = 0
$data[192,3] = 4
$data[192,4] = 7
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = "Groenlandia"
$data[193,1] = 11
$data[193,2] = 0
$data[193,3] = 11
$data[193,4] = 0
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = "Islas Turcas y Caicos"
$data[194,1] = 10
$data[194,2] = 0
$data[194,3] = 0
$data[194,4] = 9
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 1
$data[195,0] = "Surinam"
$data[195,1] = 10
$data[195,2] = 0
$data[195,3] = 6
$data[195,4] = 3
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 1
$data[196,0] = "Gambia"
$data[196,1] = 9
$data[196,2] = 0
$data[196,3] = 2
$data[196,4] = 6
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 1
$data[197,0] = "Nicaragua"
$data[197,1] = 9
$data[197,2] = 0
$data[197,3] = 4
$data[197,4] = 4
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 1
$data[198,0] = "Timor Oriental"
$data[198,1] = 8
$data[198,2] = 2
$data[198,3] = 1
$data[198,4] = 7
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = "Santa Sede"
$data[199,1] = 8
$data[199,2] = 0
$data[199,3] = 2
$data[199,4] = 6
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0
$data[200,0] = "Mauritania"
$data[200,1] = 7
$data[200,2] = 0
$data[200,3] = 2
$data[200,4] = 4
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 1
$data[201,0] = "Sahara Occidental"
$data[201,1] = 6
$data[201,2] = 0
$data[201,3] = 0
$data[201,4] = 6
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = "San Bartolome"
$data[202,1] = 6
$data[202,2] = 0
$data[202,3] = 4
$data[202,4] = 2
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = "Burundi"
$data[203,1] = 5
$data[203,2] = 0
$data[203,3] = 0
$data[203,4] = 4
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 1
$data[204,0] = "Butan"
$data[204,1] = 5
$data[204,2] = 0
$data[204,3] = 2
$data[204,4] = 3
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = "Santo Tome y Principe"
$data[205,1] = 4
$data[205,2] = 0
$data[205,3] = 0
$data[205,4] = 4
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = "Sudan del Sur"
$data[206,1] = 4
$data[206,2] = 0
$data[206,3] = 0
$data[206,4] = 4
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 0
$data[207,0] = "Bonaire, San Eustaquio y Saba"
$data[207,1] = 3
$data[207,2] = 0
$data[207,3] = 0
$data[207,4] = 3
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 0
$data[208,0] = "Anguila"
$data[208,1] = 3
$data[208,2] = 0
$data[208,3] = 1
$data[208,4] = 2
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 0
$data[209,0] = "Islas Virgenes Britanicas"
$data[209,1] = 3
$data[209,2] = 0
$data[209,3] = 2
$data[209,4] = 1
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 0
$data[210,0] = "Papua Nueva Guinea"
$data[210,1] = 2
$data[210,2] = 0
$data[210,3] = 0
$data[210,4] = 2
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 0
$data[211,0] = "Yemen"
$data[211,1] = 1
$data[211,2] = 0
$data[211,3] = 0
$data[211,4] = 1
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 0
$data[212,0] = "San Pedro y Miquelon"
$data[212,1] = 1
$data[212,2] = 0
$data[212,3] = 0
$data[212,4] = 1
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 0

$ws.Range("A4:H216").Value2 = $data
